$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2533.889
$ws.Range("I32").Value = 2210.5
$ws.Range("K32").Value = 2210.5
$ws.Range("M32").Value = -1884.5
$ws.Range("H39").Value = 1815.1666
$ws.Range("I39").Value = 1071.091
$ws.Range("K39").Value = 3213.273
$ws.Range("M39").Value = -2917.273
$ws.Range("H43").Value = 3092.5
$ws.Range("I43").Value = 2500
$ws.Range("K43").Value = 2500
$ws.Range("M43").Value = -2431
$ws.Range("H52").Value = 7484.3
$ws.Range("J52").Value = 6332.6665
$ws.Range("L52").Value = 18997.9995
$ws.Range("N52").Value = -19317.9995
$ws.Range("H53").Value = 4431.0835
$ws.Range("I53").Value = 6516
$ws.Range("J53").Value = 261.25
$ws.Range("K53").Value = 6516
$ws.Range("L53").Value = 261.25
$ws.Range("M53").Value = -5879
$ws.Range("N53").Value = -1535.25
$ws.Range("H125").Value = 1930
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1930
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 17370
$ws.Range("N125").Value = -22290
$ws.Range("H129").Value = 1507.5
$ws.Range("J129").Value = 2693.6
$ws.Range("L129").Value = 8080.799999999999
$ws.Range("N129").Value = -18080.8
$ws.Range("H132").Value = 30305954
$ws.Range("I132").Value = 35716892
$ws.Range("K132").Value = 107150676
$ws.Range("M132").Value = -107148146
$ws.Range("H138").Value = 5172.067
$ws.Range("I138").Value = 2302.276
$ws.Range("J138").Value = 6981.2827
$ws.Range("K138").Value = 6906.828
$ws.Range("L138").Value = 20943.8481
$ws.Range("M138").Value = -1766.828
$ws.Range("N138").Value = -31223.8481
$ws.Range("H141").Value = 10168.652
$ws.Range("J141").Value = 29296.6
$ws.Range("L141").Value = 87889.79999999999
$ws.Range("N141").Value = -98249.79999999999
$ws.Range("M125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5492.056
$ws.Range("I32").Value = 4275.9844
$ws.Range("K32").Value = 4275.9844
$ws.Range("M32").Value = -3988.9844
$ws.Range("H102").Value = 3529.4736
$ws.Range("I102").Value = 2222.5
$ws.Range("K102").Value = 2222.5
$ws.Range("M102").Value = -600.5
$ws.Range("H132").Value = 8501.019
$ws.Range("I132").Value = 4974.5776
$ws.Range("K132").Value = 14923.7328
$ws.Range("M132").Value = -12393.7328

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4041.796
$ws.Range("I86").Value = 3772.375
$ws.Range("J86").Value = 4548.9414
$ws.Range("K86").Value = 3772.375
$ws.Range("L86").Value = 4548.9414
$ws.Range("M86").Value = -2649.375
$ws.Range("N86").Value = -6794.9414
$ws.Range("H89").Value = 4041.796
$ws.Range("I89").Value = 3772.375
$ws.Range("J89").Value = 4548.9414
$ws.Range("K89").Value = 18861.875
$ws.Range("L89").Value = 22744.707
$ws.Range("M89").Value = -13245.875
$ws.Range("N89").Value = -33976.70699999999
$ws.Range("H94").Value = 40003492
$ws.Range("I94").Value = 2000.625
$ws.Range("K94").Value = 2000.625
$ws.Range("M94").Value = -1549.625
$ws.Range("H99").Value = 66668956
$ws.Range("I99").Value = 71430930
$ws.Range("K99").Value = 71430930
$ws.Range("M99").Value = -71429432
$ws.Range("H105").Value = 2530
$ws.Range("I105").Value = 2385.8
$ws.Range("K105").Value = 2385.8
$ws.Range("M105").Value = -638.8000000000002
$ws.Range("H134").Value = 38463390
$ws.Range("J134").Value = 4014
$ws.Range("L134").Value = 12042
$ws.Range("N134").Value = -17112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2535
$ws.Range("I16").Value = 1997
$ws.Range("J16").Value = 2700.5386
$ws.Range("K16").Value = 1997
$ws.Range("L16").Value = 2700.5386
$ws.Range("M16").Value = -1710
$ws.Range("N16").Value = -3274.5386
$ws.Range("H31").Value = 386146.28
$ws.Range("I31").Value = 7628.7393
$ws.Range("K31").Value = 7628.7393
$ws.Range("M31").Value = -7333.7393
$ws.Range("H34").Value = 386146.28
$ws.Range("I34").Value = 7628.7393
$ws.Range("K34").Value = 7628.7393
$ws.Range("M34").Value = -7426.7393
$ws.Range("H105").Value = 14423.45
$ws.Range("I105").Value = 18952.133
$ws.Range("K105").Value = 18952.133
$ws.Range("M105").Value = -17205.133
$ws.Range("H113").Value = 2535
$ws.Range("I113").Value = 1997
$ws.Range("J113").Value = 2700.5386
$ws.Range("K113").Value = 1997
$ws.Range("L113").Value = 2700.5386
$ws.Range("M113").Value = 173
$ws.Range("N113").Value = -7040.5386
$ws.Range("H124").Value = 60000
$ws.Range("J124").Value = 60000
$ws.Range("L124").Value = 60000
$ws.Range("N124").Value = -64910
$ws.Range("H132").Value = 4452.5
$ws.Range("I132").Value = 4094.5
$ws.Range("K132").Value = 12283.5
$ws.Range("M132").Value = -9753.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 13333599
$ws.Range("J2").Value = 47619796
$ws.Range("L2").Value = 285718776
$ws.Range("N2").Value = -285719002
$ws.Range("H97").Value = 511.33334
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("H98").Value = 650
$ws.Range("I98").Value = 650
$ws.Range("K98").Value = 1950
$ws.Range("M98").Value = -452
$ws.Range("H131").Value = 2000.0394
$ws.Range("I131").Value = 1900.6154
$ws.Range("J131").Value = 2020.5555
$ws.Range("K131").Value = 5701.8462
$ws.Range("L131").Value = 6061.666499999999
$ws.Range("M131").Value = -661.8462
$ws.Range("N131").Value = -16141.6665
$ws.Range("M97").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("H54").Value = 10000
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 10000
$ws.Range("N54").Value = -10780
$ws.Range("H80").Value = 2009.9615
$ws.Range("I80").Value = 1923.2142
$ws.Range("J80").Value = 2111.1667
$ws.Range("K80").Value = 1923.2142
$ws.Range("L80").Value = 2111.1667
$ws.Range("M80").Value = -925.2141999999999
$ws.Range("N80").Value = -4107.1667
$ws.Range("H83").Value = 2009.9615
$ws.Range("I83").Value = 1923.2142
$ws.Range("J83").Value = 2111.1667
$ws.Range("K83").Value = 9616.071
$ws.Range("L83").Value = 10555.8335
$ws.Range("M83").Value = -4624.071
$ws.Range("N83").Value = -20539.8335
$ws.Range("H102").Value = 3627.0833
$ws.Range("I102").Value = 3395.1667
$ws.Range("K102").Value = 3395.1667
$ws.Range("M102").Value = -1773.1667
$ws.Range("H126").Value = 5097.2085
$ws.Range("I126").Value = 4944.524
$ws.Range("K126").Value = 14833.572
$ws.Range("M126").Value = -12363.572
$ws.Range("N33").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 32262142
$ws.Range("I7").Value = 50002090
$ws.Range("K7").Value = 50002090
$ws.Range("M7").Value = -50001978
$ws.Range("H40").Value = 3747.8235
$ws.Range("I40").Value = 3794
$ws.Range("J40").Value = 3532.3333
$ws.Range("K40").Value = 3794
$ws.Range("L40").Value = 3532.3333
$ws.Range("M40").Value = -3658
$ws.Range("N40").Value = -3804.3333
$ws.Range("H61").Value = 6237.2104
$ws.Range("I61").Value = 6350.4375
$ws.Range("K61").Value = 6350.4375
$ws.Range("M61").Value = -6148.4375
$ws.Range("H105").Value = 43665.668
$ws.Range("J105").Value = 43665.668
$ws.Range("L105").Value = 43665.668
$ws.Range("N105").Value = -50653.668
$ws.Range("H113").Value = 6237.2104
$ws.Range("I113").Value = 6350.4375
$ws.Range("K113").Value = 6350.4375
$ws.Range("M113").Value = -4180.4375
$ws.Range("H122").Value = 7014.2334
$ws.Range("I122").Value = 4233.9546
$ws.Range("J122").Value = 14660
$ws.Range("K122").Value = 12701.8638
$ws.Range("L122").Value = 43980
$ws.Range("M122").Value = -10251.8638
$ws.Range("N122").Value = -48880
$ws.Range("H125").Value = 73995
$ws.Range("J125").Value = 73995
$ws.Range("L125").Value = 73995
$ws.Range("N125").Value = -83835
$ws.Range("H126").Value = 32262142
$ws.Range("I126").Value = 50002090
$ws.Range("K126").Value = 150006270
$ws.Range("M126").Value = -150003800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H126").Value = 28370478
$ws.Range("J126").Value = 55557540
$ws.Range("L126").Value = 166672620
$ws.Range("N126").Value = -166677560
$ws.Range("H132").Value = 1544.5938
$ws.Range("I132").Value = 1591.2
$ws.Range("K132").Value = 4773.6
$ws.Range("M132").Value = -2243.6
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()
